$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-7, columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T
# (F and L stay unchanged at 1)

$data = @{
    2 = @{ E=3; G=85.78116200000001; H=257.343486; I=0.6103818382160713; J=0.6103818382160713; K=3; M=28.80691533333333; N=86.42074599999999; O=0.974533567023909; P=0.974533567023909; Q=2471.090670928951; R=22239.81603836056; S=0.5948375900433185; T=0.5948375900433185 }
    3 = @{ E=3; G=85.78116200000001; H=257.343486; I=0.6103818382160713; J=0.6103818382160713; K=3; M=0.75278; N=2.25834; O=0.02546643297609089; P=0.02546643297609089; Q=64.57434313036001; R=581.1690881732401; S=0.01554424817275273; T=0.01554424817275273 }
    4 = @{ E=3; G=11.01836233333333; H=33.055087; I=0.07840192529859551; J=0.07840192529859551; K=3; M=28.80691533333333; N=86.42074599999999; O=0.974533567023909; P=0.974533567023909; Q=317.4050308483224; R=2856.645277634902; S=0.07640530792278233; T=0.07640530792278233 }
    5 = @{ E=3; G=11.01836233333333; H=33.055087; I=0.07840192529859551; J=0.07840192529859551; K=3; M=0.75278; N=2.25834; O=0.02546643297609089; P=0.02546643297609089; Q=8.294402797286667; R=74.64962517558; S=0.001996617375813167; T=0.001996617375813167 }
    6 = @{ E=3; G=43.73736033333333; H=131.212081; I=0.3112162364853332; J=0.3112162364853332; K=3; M=28.80691533333333; N=86.42074599999999; O=0.974533567023909; P=0.974533567023909; Q=1259.938436025825; R=11339.44592423243; S=0.3032906690578082; T=0.3032906690578082 }
    7 = @{ E=3; G=43.73736033333333; H=131.212081; I=0.3112162364853332; J=0.3112162364853332; K=3; M=0.75278; N=2.25834; O=0.02546643297609089; P=0.02546643297609089; Q=32.92461011172666; R=296.3214910055401; S=0.007925567427524989; T=0.007925567427524989 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
